$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values scraped on this run; D/E columns are plain text (not numbers),
# so force text format on purely-numeric-looking values to avoid Excel
# auto-converting them (which would lose exact formatting, e.g. trailing zeros).

$ws.Range('D2').Value = '62.916.18'
$ws.Range('E2').Value = '  -2.44%  '
$ws.Range('D3').Value = '3.119.18'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.61'
$ws.Range('E5').Value = '  -1.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.83'
$ws.Range('E6').Value = '  -5.35%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '3.114.82'
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.514'
$ws.Range('E9').Value = '  -1.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.144'
$ws.Range('E10').Value = '  -4.35%  '
$ws.Range('E11').Value = '  -3.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.453'
$ws.Range('E12').Value = '  -3.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000244'
$ws.Range('E13').Value = '  -5.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.96'
$ws.Range('E14').Value = '  -3.90%  '
$ws.Range('D15').Value = '3.636.09'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('E16').Value = '  +1.55%  '
$ws.Range('D17').Value = '63.010.43'
$ws.Range('E17').Value = '  -2.23%  '
$ws.Range('D18').Value = '3.122.41'
$ws.Range('E18').Value = '  -1.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.63'
$ws.Range('E19').Value = '  -3.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '469.51'
$ws.Range('E20').Value = '  -2.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.05'
$ws.Range('E21').Value = '  -3.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.693'
$ws.Range('E22').Value = '  -2.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.63'
$ws.Range('E23').Value = '  -1.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.47'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.88'
$ws.Range('E25').Value = '  -4.16%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.71'
$ws.Range('E27').Value = '  -1.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.81'
$ws.Range('E28').Value = '  -7.14%  '
$ws.Range('E29').Value = '  +1.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.80'
$ws.Range('E30').Value = '  -6.16%  '
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.53'
$ws.Range('E32').Value = '  -1.68%  '
$ws.Range('E33').Value = '  -5.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.51'
$ws.Range('E34').Value = '  -5.75%  '
$ws.Range('E35').Value = '  -3.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.74'
$ws.Range('E36').Value = '  -4.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.90'
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('D38').Value = '0.0₃0689'
$ws.Range('E38').Value = '  -11.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0385'
$ws.Range('E39').Value = '  -2.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '415.36'
$ws.Range('E40').Value = '  -6.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.18'
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('D42').Value = '2.894.29'
$ws.Range('E42').Value = '  +1.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.64'
$ws.Range('E43').Value = '  -12.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.112'
$ws.Range('E44').Value = '  -5.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.260'
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('E47').Value = '  -6.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.26'
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.23'
$ws.Range('E50').Value = '  -7.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.19'
$ws.Range('E51').Value = '  +0.06%  '
